{"js": "// Bottom command bar for Windows Phone\n// 1) Insert a new bullet paragraph (\"The band tile will notify you...\")\n//    right before the \"Using the latest SDK...\" paragraph.\n// 2) Turn the old \"- When opening the band tile...\" paragraph into\n//    \"- Note: when opening the band tile...\" split across three runs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their (pre-edit) text.\nlet sdkParagraph = null;\nlet noteParagraph = null;\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"Using the latest SDK\") !== -1) {\n    sdkParagraph = p;\n  } else if (t.indexOf(\"When opening the band tile\") !== -1) {\n    noteParagraph = p;\n  }\n}\n\nif (!sdkParagraph || !noteParagraph) {\n  throw new Error(\"Could not locate expected paragraphs in the document.\");\n}\n\n// 1) Insert the new paragraph just before the \"Using the latest SDK...\" one.\nsdkParagraph.insertParagraph(\n  \"- The band tile will notify you when your paired Windows phone reaches 100% charge state.\",\n  Word.InsertLocation.before\n);\n\n// 2) Replace the \"- When opening...\" paragraph's content with three runs:\n//      \"- \"  +  \"Note: w\"  +  \"hen opening the band tile, please allow a few seconds for the tile to update.\"\nconst noteOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">- </w:t></w:r>' +\n  '<w:r><w:t>Note: w</w:t></w:r>' +\n  '<w:r><w:t>hen opening the band tile, please allow a few seconds for the tile to update.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nnoteParagraph.insertOoxml(noteOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Bottom command bar for Windows Phone\n#\n# 1) Insert a new bullet paragraph (\"The band tile will notify you...\")\n#    right before the \"Using the latest SDK...\" paragraph.\n# 2) Turn the old \"- When opening the band tile...\" paragraph into\n#    \"- Note: when opening the band tile...\" split across three runs.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the new paragraph ------------------------------------\n\n$sdkPara = $null\n$paraBeforeSdk = $null\n$prev = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Using the latest SDK*\") {\n        $sdkPara = $p\n        $paraBeforeSdk = $prev\n        break\n    }\n    $prev = $p\n}\n\nif ($sdkPara -eq $null) {\n    throw \"Could not locate the 'Using the latest SDK' paragraph.\"\n}\n\n# Create a blank paragraph immediately before it, then fill in its text.\n$sdkPara.Range.InsertParagraphBefore()\n\n$newPara = $null\nif ($paraBeforeSdk -eq $null) {\n    $newPara = $d.Paragraphs.First\n} else {\n    $newPara = $paraBeforeSdk.Next()\n}\n$newPara.Range.Text = \"- The band tile will notify you when your paired Windows phone reaches 100% charge state.\"\n\n# --- Step 2: rewrite the \"When opening...\" paragraph into 3 runs ---------\n\n$notePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*When opening the band tile*\") {\n        $notePara = $p\n        break\n    }\n}\n\nif ($notePara -eq $null) {\n    throw \"Could not locate the 'When opening the band tile' paragraph.\"\n}\n\n$ooxmlFragment = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:r><w:t xml:space=\"preserve\">- </w:t></w:r>' +\n    '<w:r><w:t>Note: w</w:t></w:r>' +\n    '<w:r><w:t>hen opening the band tile, please allow a few seconds for the tile to update.</w:t></w:r>' +\n    '</w:p>'\n\n$notePara.Range.InsertXML($ooxmlFragment)\n"}
